$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update existing rows: Yes -> No (row 8, Profile) and Yes -> NO (row 9, UserDashboard)
$ws.Range("C8").Value = "No"
$ws.Range("C9").Value = "NO"

# Add new row 10 for "Admin Tool section-Invoice Lookup" test case (Tools)
$ws.Range("A10").Value = "Tools"
$ws.Range("B10").Value = "Admin Tools Scenarios"
$ws.Range("C10").Value = "Yes"

# Update selection to match the target state
$ws.Range("B14").Select()
